# Changes of 28th july 2022
# Update row 5 tracking number (ShipmentTracking) and actual rate (ActualRate)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces these numeric-looking values to be stored as
# text (matching the original cell type), then resetting the Style back
# to "Normal" clears the auto-added quote-prefix / number-format style so
# the cell keeps using the default (unstyled) format, same as before.
$ws.Range("P5").Value = "'320018767480"
$ws.Range("P5").Style = "Normal"

$ws.Range("Q5").Value = "'`$46.27"
$ws.Range("Q5").Style = "Normal"
